# Apply weekly-refreshed Fruta/Hortaliza price data (Maracuyá, Vega Modelo de Temuco)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44726
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 34000
$ws.Range("O2").Value = 34000
$ws.Range("P2").Value = 34000
$ws.Range("S2").Value = 1889
$ws.Range("D3").Value = 44392
$ws.Range("M3").Value = 20
$ws.Range("N3").Value = 35000
$ws.Range("O3").Value = 35000
$ws.Range("P3").Value = 35000
$ws.Range("S3").Value = 1944
$ws.Range("D4").Value = 44442
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = 35000
$ws.Range("O4").Value = 35000
$ws.Range("P4").Value = 35000
$ws.Range("R4").Value = 'Perú'
$ws.Range("S4").Value = 1944
$ws.Range("D5").Value = 44438
$ws.Range("M5").Value = 25
$ws.Range("R5").Value = 'Región de Arica y Parinacota'
$ws.Range("D6").Value = 44363
$ws.Range("M6").Value = 144
$ws.Range("N6").Value = 1700
$ws.Range("O6").Value = 1700
$ws.Range("P6").Value = 1700
$ws.Range("Q6").Value = '$/kilo'
$ws.Range("S6").Value = 1700
$ws.Range("T6").Value = 1
$ws.Range("D7").Value = 44357
$ws.Range("M7").Value = 10
$ws.Range("N7").Value = 38000
$ws.Range("O7").Value = 38000
$ws.Range("P7").Value = 38000
$ws.Range("R7").Value = 'Perú'
$ws.Range("S7").Value = 2111
$ws.Range("D8").Value = 44676
$ws.Range("M8").Value = 55
$ws.Range("N8").Value = 28000
$ws.Range("O8").Value = 30000
$ws.Range("P8").Value = 28909
$ws.Range("S8").Value = 1606
$ws.Range("D9").Value = 44377
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 40000
$ws.Range("O9").Value = 40000
$ws.Range("P9").Value = 40000
$ws.Range("S9").Value = 2222
$ws.Range("D10").Value = 44662
$ws.Range("M10").Value = 15
$ws.Range("N10").Value = 30000
$ws.Range("O10").Value = 30000
$ws.Range("P10").Value = 30000
$ws.Range("S10").Value = 1667
$ws.Range("D11").Value = 44433
$ws.Range("M11").Value = 15
$ws.Range("R11").Value = 'Región de Arica y Parinacota'
$ws.Range("D13").Value = 44679
$ws.Range("M13").Value = 35
$ws.Range("N13").Value = 34000
$ws.Range("O13").Value = 34000
$ws.Range("P13").Value = 34000
$ws.Range("S13").Value = 1889
$ws.Range("D14").Value = 44679
$ws.Range("M14").Value = 55
$ws.Range("N14").Value = 28000
$ws.Range("O14").Value = 28000
$ws.Range("P14").Value = 28000
$ws.Range("S14").Value = 1556
$ws.Range("D15").Value = 44279
$ws.Range("M15").Value = 30
$ws.Range("O15").Value = 36000
$ws.Range("P15").Value = 35667
$ws.Range("S15").Value = 1982
$ws.Range("D16").Value = 44719
$ws.Range("M16").Value = 25
$ws.Range("N16").Value = 34000
$ws.Range("O16").Value = 34000
$ws.Range("P16").Value = 34000
$ws.Range("R16").Value = 'Región de Arica y Parinacota'
$ws.Range("S16").Value = 1889
$ws.Range("D17").Value = 44424
$ws.Range("M17").Value = 15
$ws.Range("N17").Value = 35000
$ws.Range("O17").Value = 35000
$ws.Range("P17").Value = 35000
$ws.Range("S17").Value = 1944
$ws.Range("D18").Value = 44431
$ws.Range("M18").Value = 30
$ws.Range("N18").Value = 35000
$ws.Range("O18").Value = 35000
$ws.Range("P18").Value = 35000
$ws.Range("S18").Value = 1944
$ws.Range("D19").Value = 44405
$ws.Range("M19").Value = 10
$ws.Range("N19").Value = 35000
$ws.Range("O19").Value = 35000
$ws.Range("P19").Value = 35000
$ws.Range("S19").Value = 1944
$ws.Range("D20").Value = 44434
$ws.Range("M20").Value = 40
$ws.Range("N20").Value = 35000
$ws.Range("O20").Value = 35000
$ws.Range("P20").Value = 35000
$ws.Range("S20").Value = 1944
$ws.Range("D21").Value = 44379
$ws.Range("M21").Value = 10
$ws.Range("N21").Value = 30000
$ws.Range("O21").Value = 30000
$ws.Range("P21").Value = 30000
$ws.Range("S21").Value = 1667
$ws.Range("D22").Value = 44432
$ws.Range("M22").Value = 10
$ws.Range("R22").Value = 'Perú'
$ws.Range("D23").Value = 44671
$ws.Range("M23").Value = 20
$ws.Range("N23").Value = 32000
$ws.Range("O23").Value = 32000
$ws.Range("P23").Value = 32000
$ws.Range("Q23").Value = '$/caja 18 kilos'
$ws.Range("S23").Value = 1778
$ws.Range("T23").Value = 18
$ws.Range("D24").Value = 44658
$ws.Range("M24").Value = 30
$ws.Range("N24").Value = 28000
$ws.Range("O24").Value = 28000
$ws.Range("P24").Value = 28000
$ws.Range("S24").Value = 1556
$ws.Range("D25").Value = 44448
$ws.Range("M25").Value = 50
$ws.Range("N25").Value = 38000
$ws.Range("O25").Value = 38000
$ws.Range("P25").Value = 38000
$ws.Range("R25").Value = 'Región de Arica y Parinacota'
$ws.Range("S25").Value = 2111
$ws.Range("D26").Value = 44645
$ws.Range("M26").Value = 5
$ws.Range("N26").Value = 30000
$ws.Range("O26").Value = 30000
$ws.Range("P26").Value = 30000
$ws.Range("S26").Value = 1667
$ws.Range("D27").Value = 44721
$ws.Range("M27").Value = 5
$ws.Range("N27").Value = 35000
$ws.Range("O27").Value = 35000
$ws.Range("P27").Value = 35000
$ws.Range("R27").Value = 'Perú'
$ws.Range("S27").Value = 1944
$ws.Range("D28").Value = 44264
$ws.Range("M28").Value = 20
$ws.Range("N28").Value = 40000
$ws.Range("O28").Value = 40000
$ws.Range("P28").Value = 40000
$ws.Range("R28").Value = 'Región de Arica y Parinacota'
$ws.Range("S28").Value = 2222
$ws.Range("D29").Value = 44699
$ws.Range("M29").Value = 20
$ws.Range("D30").Value = 44664
$ws.Range("M30").Value = 15
$ws.Range("N30").Value = 30000
$ws.Range("O30").Value = 30000
$ws.Range("P30").Value = 30000
$ws.Range("R30").Value = 'Perú'
$ws.Range("S30").Value = 1667
$ws.Range("D31").Value = 44690
$ws.Range("M31").Value = 25
$ws.Range("N31").Value = 34000
$ws.Range("O31").Value = 34000
$ws.Range("P31").Value = 34000
$ws.Range("S31").Value = 1889
$ws.Range("D32").Value = 44435
$ws.Range("M32").Value = 10
$ws.Range("N32").Value = 35000
$ws.Range("O32").Value = 35000
$ws.Range("P32").Value = 35000
$ws.Range("S32").Value = 1944
$ws.Range("D33").Value = 44435
$ws.Range("M33").Value = 105
$ws.Range("N33").Value = 35000
$ws.Range("O33").Value = 35000
$ws.Range("P33").Value = 35000
$ws.Range("S33").Value = 1944
$ws.Range("D34").Value = 44294
$ws.Range("D35").Value = 44629
$ws.Range("M35").Value = 20
$ws.Range("N35").Value = 35000
$ws.Range("O35").Value = 35000
$ws.Range("P35").Value = 35000
$ws.Range("R35").Value = 'Región de Arica y Parinacota'
$ws.Range("S35").Value = 1944
$ws.Range("D36").Value = 44418
$ws.Range("M36").Value = 30
$ws.Range("D37").Value = 44704
$ws.Range("M37").Value = 25
$ws.Range("N37").Value = 35000
$ws.Range("O37").Value = 35000
$ws.Range("P37").Value = 35000
$ws.Range("S37").Value = 1944
$ws.Range("D38").Value = 44669
$ws.Range("M38").Value = 40
$ws.Range("N38").Value = 32000
$ws.Range("O38").Value = 32000
$ws.Range("P38").Value = 32000
$ws.Range("S38").Value = 1778
$ws.Range("D39").Value = 44449
$ws.Range("N39").Value = 38000
$ws.Range("O39").Value = 38000
$ws.Range("P39").Value = 38000
$ws.Range("S39").Value = 2111
$ws.Range("D40").Value = 44369
$ws.Range("M40").Value = 5
$ws.Range("N40").Value = 35000
$ws.Range("O40").Value = 35000
$ws.Range("P40").Value = 35000
$ws.Range("Q40").Value = '$/caja 18 kilos'
$ws.Range("R40").Value = 'Perú'
$ws.Range("S40").Value = 1944
$ws.Range("T40").Value = 18
$ws.Range("D41").Value = 44364
$ws.Range("M41").Value = 90
$ws.Range("N41").Value = 1700
$ws.Range("O41").Value = 1700
$ws.Range("P41").Value = 1700
$ws.Range("Q41").Value = '$/kilo'
$ws.Range("S41").Value = 1700
$ws.Range("T41").Value = 1
$ws.Range("D42").Value = 44634
$ws.Range("N42").Value = 45000
$ws.Range("O42").Value = 45000
$ws.Range("P42").Value = 45000
$ws.Range("S42").Value = 2500
